$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.719.52'
$ws.Range('E2').Value = '  +2.88%  '
$ws.Range('D3').Value = '2.237.41'
$ws.Range('E3').Value = '  +0.75%  '
$ws.Range('E4').Value = '  -0.10%  '
$__style = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '232.41'
$ws.Range('D5').Style = $__style
$ws.Range('E5').Value = '  +0.64%  '
$__style = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.622'
$ws.Range('D6').Style = $__style
$ws.Range('E6').Value = '  -0.55%  '
$__style = $ws.Range('D7').Style
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '61.14'
$ws.Range('D7').Style = $__style
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('E8').Value = '  +0.11%  '
$__style = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.407'
$ws.Range('D9').Style = $__style
$ws.Range('E9').Value = '  +1.04%  '
$__style = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0916'
$ws.Range('D10').Style = $__style
$ws.Range('E10').Value = '  +3.13%  '
$ws.Range('E11').Value = '  +0.11%  '
$ws.Range('D12').Value = '2.569.67'
$ws.Range('E12').Value = '  +0.80%  '
$__style = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '15.60'
$ws.Range('D13').Style = $__style
$ws.Range('E13').Value = '  -0.32%  '
$__style = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '22.37'
$ws.Range('D14').Style = $__style
$ws.Range('E14').Value = '  +2.39%  '
$__style = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.68'
$ws.Range('D15').Style = $__style
$ws.Range('E15').Value = '  +2.25%  '
$__style = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.802'
$ws.Range('D16').Style = $__style
$ws.Range('E16').Value = '  +0.43%  '
$ws.Range('D17').Value = '2.238.72'
$ws.Range('E17').Value = '  +0.89%  '
$ws.Range('D18').Value = '42.516.07'
$ws.Range('E18').Value = '  +2.88%  '
$ws.Range('D19').Value = '0.0₃0945'
$ws.Range('E19').Value = '  +5.35%  '
$__style = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.20'
$ws.Range('D20').Style = $__style
$ws.Range('E20').Value = '  +2.38%  '
$__style = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.51'
$ws.Range('D21').Style = $__style
$ws.Range('E21').Value = '  -0.61%  '
$__style = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '245.38'
$ws.Range('D22').Style = $__style
$ws.Range('E22').Value = '  -2.35%  '
$ws.Range('B23').Value = 'PancakeSwap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$__style = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.62'
$ws.Range('D23').Style = $__style
$ws.Range('E23').Value = '  +9.31%  '
$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$__style = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.999'
$ws.Range('D24').Style = $__style
$ws.Range('E24').Value = '  -0.13%  '
$__style = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.40'
$ws.Range('D25').Style = $__style
$ws.Range('E25').Value = '  +5.43%  '
$__style = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.75'
$ws.Range('D26').Style = $__style
$ws.Range('E26').Value = '  +1.92%  '
$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$__style = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.143'
$ws.Range('D27').Style = $__style
$ws.Range('E27').Value = '  +2.40%  '
$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$__style = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '169.59'
$ws.Range('D28').Style = $__style
$ws.Range('E28').Value = '  +0.82%  '
$__style = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '20.42'
$ws.Range('D29').Style = $__style
$ws.Range('E29').Value = '  +2.27%  '
$ws.Range('E30').Value = '  +3.05%  '
$ws.Range('E31').Value = '  +1.07%  '
$ws.Range('E32').Value = '  -1.25%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$__style = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.01'
$ws.Range('D33').Style = $__style
$ws.Range('E33').Value = '  +0.70%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$__style = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.74'
$ws.Range('D34').Style = $__style
$ws.Range('E34').Value = '  +2.31%  '
$__style = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0656'
$ws.Range('D35').Style = $__style
$ws.Range('E35').Value = '  +4.74%  '
$__style = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.43'
$ws.Range('D36').Style = $__style
$ws.Range('E36').Value = '  -2.44%  '
$__style = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.38'
$ws.Range('D37').Style = $__style
$ws.Range('E37').Value = '  +0.98%  '
$__style = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.60'
$ws.Range('D38').Style = $__style
$ws.Range('E38').Value = '  -2.64%  '
$__style = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0250'
$ws.Range('D39').Style = $__style
$ws.Range('E39').Value = '  +4.64%  '
$ws.Range('E40').Value = '  -0.10%  '
$__style = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.68'
$ws.Range('D41').Style = $__style
$ws.Range('E41').Value = '  +0.95%  '
$__style = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.000227'
$ws.Range('D42').Style = $__style
$ws.Range('E42').Value = '  -7.46%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$__style = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.21'
$ws.Range('D43').Style = $__style
$ws.Range('E43').Value = '  +1.30%  '
$ws.Range('B44').Value = 'Cronos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$__style = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0965'
$ws.Range('D44').Style = $__style
$ws.Range('E44').Value = '  -2.14%  '
$__style = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '97.34'
$ws.Range('D45').Style = $__style
$ws.Range('E45').Value = '  -1.74%  '
$__style = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.40'
$ws.Range('D46').Style = $__style
$ws.Range('E46').Value = '  -8.95%  '
$ws.Range('D47').Value = '1.459.32'
$ws.Range('E47').Value = '  -0.42%  '
$__style = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '16.31'
$ws.Range('D48').Style = $__style
$ws.Range('E48').Value = '  -1.76%  '
$__style = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.75'
$ws.Range('D49').Style = $__style
$ws.Range('E49').Value = '  -1.13%  '
$ws.Range('E50').Value = '  -0.47%  '
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$__style = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.23'
$ws.Range('D51').Style = $__style
$ws.Range('E51').Value = '  +4.15%  '
